$d = $word.ActiveDocument

# 1. Template city: Екатеринбург -> Новороссийск (the date-prefix spacing run
#    merges automatically with the city run since formatting is identical).
[void]$d.Content.Find.Execute("Екатеринбург", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Новороссийск", 2)

# 2. Crypto num tag: tighten the Normal style's auto line spacing
#    from 266 (13.3pt) to 264 (13.2pt) twentieths.
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.LineSpacingRule = 5
$normalStyle.ParagraphFormat.LineSpacing = 13.2
